$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: new time entry (date / task / hours) ---
$ws.Range("D17").Value = 45430
$ws.Range("E17").Value = "Ultimas modificaciones del proyecto"
$ws.Range("F17").Value = 5

# --- Row 19: move the "MEMORIA" section label here and bold it ---
$ws.Range("B19").Value = "MEMORIA"
$ws.Range("B19").Font.Bold = $true

# --- Row 20: new time entry ---
$ws.Range("D20").Value = 45430
$ws.Range("E20").Value = "Avances en multiples apartados"
$ws.Range("F20").Value = 3

# --- Row 21: remove the old "MEMORIA" label, add a new time entry ---
$ws.Range("B21").ClearContents()
$ws.Range("D21").Value = 45431
$ws.Range("E21").Value = "Ultimas modificaciones"
$ws.Range("F21").Value = 2

# --- Row 22 becomes the TOTAL HORAS / subtotal row (was row 24) ---
$ws.Range("B22").ClearContents()
$ws.Range("D22").Value = "TOTAL HORAS"
$ws.Range("E22").ClearContents()
$ws.Range("F22").Formula = "=SUBTOTAL(109,Tabla1[Horas])"

# --- Rows 23 & 24 become blank padding rows (matching rows below them) ---
$ws.Range("B23").ClearContents()
$ws.Range("D23:F23").ClearContents()
$ws.Range("D23:F23").Style = "Normal"

$ws.Range("D24:F24").ClearContents()
$ws.Range("D24:F24").Style = "Normal"

# --- Table shrinks back to the real data range ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("D3:F21"))

# --- Drop the now-unused trailing blank rows at the bottom of the sheet ---
$ws.Rows("46:47").Delete()

# --- Update selection / active cell to match the latest edit location ---
$ws.Range("F12").Select()
